# Update the "取得日時" (retrieved datetime) timestamps in column A
# from "2025-10-17 18:22:21" to "2025-10-17 18:29:50" for rows 2-11
# on the "ランサーズ" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-17 18:29:50"

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
